# ------------------------------------------------------------------
# Applies two changes to the document:
#  1. Remove the stray "_GoBack" bookmark that sits after "... wie ein
#     funktionserweitertes Interface".
#  2. Split the "Von Überladung ..." paragraph into several runs,
#     inserting two new sentences, and re-add the "_GoBack" bookmark
#     in its new spot (right after the parenthetical remark).
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Change 1: drop the old _GoBack bookmark -----------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Change 2: rewrite the "Von Überladung ..." paragraph ----------
$oldText = "Von Überladung sprechen wir, wenn der Aufruf einer Operation anhand des konkreten Typs von Variablen oder Konstanten auf eine Methode abgebildet wird. Im Gegensatz zur dynamischen Polymorphie spielen die Inhalte der Variablen bei der Entscheidung, welche konkrete Methode aufgerufen wird, keine Rolle. Überladung kann nur von Sprachen mit statischem Typsystem unterstützt werden."

$rng = $d.Content
$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Remove the original single-run sentence …
    $rng.Delete()

    # … and rebuild it as several runs (plus the relocated bookmark)
    # by inserting a raw OOXML fragment at the same spot. This mirrors
    # how Word keeps distinct runs around an inline bookmark instead
    # of collapsing everything back into one run.
    $insertRange = $d.Range($rng.Start, $rng.Start)

    $xml = @"
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:body>
<w:p>
<w:r><w:t xml:space='preserve'>Von Überladung sprechen wir, wenn der Aufruf einer Operation anhand des konkreten Typs von Variablen oder Konstanten auf eine Methode abgebildet wird. Im Gegensatz zur dynamischen Polymorphie spielen die Inhalte der Variablen bei der Entscheidung, welche konkrete Methode aufgerufen wird, keine Rolle. </w:t></w:r>
<w:r><w:t xml:space='preserve'>So kann der Methodenname gleich sein, wodurch nur nach Parametern </w:t></w:r>
<w:r><w:t xml:space='preserve'>(unterschiedliche Anzahl, unterschiedliche Typen) </w:t></w:r>
<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>
<w:r><w:t>unterschieden wird.</w:t></w:r>
<w:r><w:t xml:space='preserve'> </w:t></w:r>
<w:r><w:t>Überladung kann nur von Sprachen mit statischem Typsystem unterstützt werden.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

    $insertRange.InsertXML($xml)
}
